$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a few timing values (rows 5-7) ---
$ws.Range("B5").Value = 0.0005550384521484375
$ws.Range("B6").Value = 0.0005118846893310547
$ws.Range("B7").Value = 0.001688718795776367

# --- Change tuple-looking text to list-looking text ---
$ws.Range("A8").Value = "[[2, 2], [2, 0], [0, 2], [0, 0], [2, 1], [1, 2], [1, 0], [0, 1], [1, 1]]"
$ws.Range("A47").Value = "[[2, 1], [1, 0], [0, 2], [0, 0], [1, 1], [0, 1], [1, 2], [2, 0], [2, 2]]"
$ws.Range("A87").Value = "[[0, 2], [0, 1], [0, 0], [1, 1], [2, 1], [1, 0], [2, 0], [1, 2], [2, 2]]"

# --- Insert a new row above the old row 110 ("Movement times") ---
# This shifts the old rows 110-114 down to 111-115.
$ws.Rows.Item(110).Insert()

# --- Populate the newly inserted row 110 with the move_fidelity entry ---
$ws.Range("A110").Value = "move_fidelity"
$ws.Range("B110").Value = 0.9992305434331438

# --- Update the "total time:" value, now located at row 114 ---
$ws.Range("B114").Value = 0.01288104057312012
